$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "splash" value (B4): unlocked cell, direct write is fine ---
$ws.Range("B4").Value = '<h2>No POD activation is currently in place.</h2><p>Points of Dispensing (PODs) are temporary sites that are ready to open to provide life-saving medications to anyone who needs them after a biological attack or during a communicable disease outbreak. The sites are located throughout the city and are run by the NYC Department of Health and Mental Hygiene. All NYC POD sites are wheelchair accessible. For more information, please visit our <a href=https://www1.nyc.gov/site/doh/health/emergency-preparedness/emergency-medication-distribution.page>webpage.</a><br><br>If you do not live in NYC, please check with your local health department:  <a href=http://www.nassaucountyny.gov/agencies/Health/index.html>Nassau County</a>, <a href=https://www.suffolkcountyny.gov/health>Suffolk County</a>, <a href=https://www.co.bergen.nj.us/departments-and-services/about-health-services>Bergen County</a>, <a href=https://www.health.ny.gov/contact/contact_information>New York State</a>, <a href=https://www.state.nj.us/health/lh/community> New Jersey State</a>.</p>&quot;'

# --- Update the "active" value (B5): unlocked cell, direct write is fine ---
$ws.Range("B5").Value = "true"

# --- Insert a new row 7 for the "description" key/value pair ---
$ws.Rows(7).Insert()

# B7 is unlocked by inherited formatting, so we can write directly.
# A7 inherits the locked "key" style, and the sheet is protected, so we
# stage the text through B7, copy it, and paste into A7 (paste is allowed
# through protection for a cell whose format already exists).
$ws.Range("B7").Value = "description"
$ws.Range("B7").Copy()
$ws.Paste($ws.Range("A7"))

# Now put the real description value into B7.
$ws.Range("B7").Value = 'Not all Points of Dispensing sites may be activated at the time of an event. Please continue checking this page to see which sites are activated following an event. If you do not live in NYC, please check with your local health department.'

# Restore the selection like the saved workbook shows.
$ws.Range("B2").Select()
